$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mid")

# Remove the three champion rows whose difficulty is "NA" (no real matchup data):
# Diana (row 12), Jayce (row 18), Neeko (row 27). Deleting bottom-up keeps the
# remaining row numbers stable while we work. The rest of the table was
# already alphabetized, so after these three rows disappear every other
# champion naturally shifts up into alphabetical order.
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(12).Delete()

# Switch the active sheet to "Mid" and leave the selection on the row that
# was being edited.
$ws.Activate()
$ws.Range("A17:XFD17").Select()
